$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column F (rows 4-32) with the new light-green color (FF92D050),
# which creates the new fill/cellXf pair used by the target style.
$ws.Range("F4:F32").Interior.Color = 5296274

# Any previously-blank cell in F4:F32 becomes an explicit 0.
for ($r = 4; $r -le 32; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $v = $cell.Value()
    if ([string]::IsNullOrEmpty($v)) {
        $cell.Value = 0
    }
}

# G6 and G21 pick up the same new fill (values already 5, unchanged).
$ws.Range("G6").Interior.Color = 5296274
$ws.Range("G21").Interior.Color = 5296274

# Column E rows 26-32 also pick up the new fill (values unchanged).
$ws.Range("E26:E32").Interior.Color = 5296274

# Move the active selection from H8 to G6.
$ws.Range("G6").Select() | Out-Null
